# TC32_Verify_store_location.xlsx — "Final changes for LogixalQA"
#
# The StoreLocator verification steps (rows 17-23 on the test-case sheet,
# plus their matching object rows on the Testdata sheet) are replaced with a
# new "Contact Us" form flow: open the Contact Us menu, pick a reason from
# the dropdown, type a message, submit it and verify the confirmation text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # TC32_Verify_store_location
$ws2 = $wb.Worksheets.Item(2)   # Testdata

# ---------------------------------------------------------------------
# Sheet 1: TC32_Verify_store_location
# ---------------------------------------------------------------------

# Row 17: was "VERIFY_WEBELEMENT_PRESENT / ContactList" -> now verifies the
# Contact Us header is present.
$ws1.Range("C17").Value = "ContactUSHeader"
$ws1.Range("E17").Value = "ContactUSHeader"

# Row 18: was "VERIFY_WEBELEMENT_PRESENT / Storelocator" -> now selects a
# reason from the Contact Us dropdown.
$ws1.Range("B18").Value = "DROPDOWN_SELECT"
$ws1.Range("C18").Value = "ContactUsDropdown"
$ws1.Range("E18").Value = "ContactUsDropdown"

# Row 19: was "ENTERTEXT / StorelocatorZip" -> now clicks into the message
# box before typing (no data-descriptor needed for this step).
$ws1.Range("B19").Value = "CLICK_PRE_ENTERTEXT"
$ws1.Range("C19").Value = "ContactUsMessage"
$ws1.Range("E19").Value = ""

# Row 20: was "CLICK / StorelocatorSearchButton" -> now enters the message
# text into the Contact Us message box.
$ws1.Range("B20").Value = "ENTERTEXT"
$ws1.Range("C20").Value = "ContactUsMessage"
$ws1.Range("E20").Value = "ContactUsMessage"

# Row 21: was "VERIFY_WEBELEMENT_PRESENT / ContactList" -> now clicks the
# button that sends the contact-us details.
$ws1.Range("B21").Value = "CLICK"
$ws1.Range("C21").Value = "SendContactusDetails"
$ws1.Range("E21").Value = ""

# Row 22: was "VERIFY_TEXT_PRESENT / StoreHeader" -> now verifies the
# confirmation element is present.
$ws1.Range("B22").Value = "VERIFY_WEBELEMENT_PRESENT"
$ws1.Range("C22").Value = "ContactUsConfirmation"
$ws1.Range("E22").Value = "ContactUsConfirmation"

# Row 23: was "CLICK / MyaccountSection" -> now verifies the confirmation
# text itself.
$ws1.Range("B23").Value = "VERIFY_TEXT_PRESENT"
$ws1.Range("C23").Value = "ContactUsConfirmation"
$ws1.Range("E23").Value = "Confirmation"

# Row 24: shift the old "CLICK / Logout" row's object down - it's now
# preceded by the My Account click that used to live on row 23.
$ws1.Range("C24").Value = "MyaccountSection"

# Row 25 (new row): CLICK / Logout, completing the flow again at the end.
$ws1.Range("A25").Value = ""
$ws1.Range("B25").Value = "CLICK"
$ws1.Range("C25").Value = "Logout"
$ws1.Range("D25").Value = "CSS"
$ws1.Range("E25").Value = ""
$ws1.Range("A25:E25").Borders.LineStyle = 1

# Column C is now wider to fit the longer "ContactUsConfirmation" values.
$ws1.Columns.Item(3).ColumnWidth = 26.1666666666667

$ws1.Activate()
$ws1.Range("E22").Select()

# ---------------------------------------------------------------------
# Sheet 2: Testdata
# ---------------------------------------------------------------------

# Row 5: object renamed from ContactList to ContactUSHeader (value stays TRUE).
$ws2.Range("A5").Value = "ContactUSHeader"

# New rows describing the Contact Us objects/data used above.
$ws2.Range("A11").Value = "ContactUsDropdown"
$ws2.Range("B11").Value = "Other"

$ws2.Range("A12").Value = "ContactUsMessage"
$ws2.Range("B12").Value = "Test Message created via Test Automation"

$ws2.Range("A13").Value = "ContactUsConfirmation"
$ws2.Range("B13").Value = $true

$ws2.Range("A14").Value = "Confirmation"
$ws2.Range("B14").Value = "Your message has been sent."

# Match the bordered look of the rest of the object/value table.
$ws2.Range("A11:B14").Borders.LineStyle = 1

# Column B is now wider to fit the longer confirmation text.
$ws2.Columns.Item(2).ColumnWidth = 38.25

$ws2.Activate()
$ws2.Range("A11:B14").Select()

$ws1.Activate()
